$d = $word.ActiveDocument

$xml6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>a.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">Name of Nominee (as it will appear on the SF 1402, Certificate of Contracting Officer Appointment) </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(6).Range.InsertXML($xml6)

$xml7 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>b.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Grade/Series or Military Rank</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(7).Range.InsertXML($xml7)

$xml8 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>c.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Position/Title</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(8).Range.InsertXML($xml8)

$xml9 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>d.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Description of Duties</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(9).Range.InsertXML($xml9)

$xml10 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>e.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Activity Name/DODAAC</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(10).Range.InsertXML($xml10)

$xml14 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>a.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Education*</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(14).Range.InsertXML($xml14)

$xml15 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>b.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Specialized Experience</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(15).Range.InsertXML($xml15)

$xml16 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>c.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Relevant Training*</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(16).Range.InsertXML($xml16)

$xml17 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>d.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Other Qualifications (</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>e.g.</w:t></w:r><w:r><w:t xml:space="preserve"> DAWIA Career Field Certifications, Defense Acquisition Corps Membership)* </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(17).Range.InsertXML($xml17)

$xml23 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>a.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Type of Appointment (</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>e.g</w:t></w:r><w:r><w:t>. PCO, ACO, Ordering Officer)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(23).Range.InsertXML($xml23)

$xml24 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>b.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Dollar Threshold (</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>e.g</w:t></w:r><w:r><w:t xml:space="preserve">. Unlimited, Limited to Specific Dollar Value) </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(24).Range.InsertXML($xml24)

$xml25 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>c.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Type of work/products/services nominee will be authorized to procure or administer</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(25).Range.InsertXML($xml25)

$xml26 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>d.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Limitations (e.g. Contract Type; Competitive or Non-Competitive Actions, etc.)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(26).Range.InsertXML($xml26)

$xml27 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2295"/></w:tabs><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>e.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">Appointment Term </w:t></w:r><w:r><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(27).Range.InsertXML($xml27)

$xml31 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>a.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Explain why this delegation is necessary and how it results in a more efficient execution and administration of the HCAs contracting operations.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(31).Range.InsertXML($xml31)

$xml32 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>b.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">Identify the anticipated workload to support the type of appointment, including the available resources at the Nominee’s Organization. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(32).Range.InsertXML($xml32)

$xml33 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>c.</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:t>Describe the internal controls that are or will be in place at the Nominee’s Organization and at the HCA/delegating office to ensure efficient and effective execution and management oversight of delegated authority</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(33).Range.InsertXML($xml33)
